$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M96").Value = -211.2307000000001
$ws.Range("H96").Value = 444.16666
$ws.Range("I96").Value = 528.0769
$ws.Range("K96").Value = 1584.2307
$ws.Range("L97").Value = 3119.7
$ws.Range("J97").Value = 1039.9
$ws.Range("N97").Value = -4111.700000000001
$ws.Range("H97").Value = 1039.9
$ws.Range("H98").Value = 739.1539
$ws.Range("I98").Value = 813.3182
$ws.Range("K98").Value = 813.3182
$ws.Range("L98").Value = 331.25
$ws.Range("N98").Value = -3327.25
$ws.Range("J98").Value = 331.25
$ws.Range("M98").Value = 684.6818
$ws.Range("K99").Value = 733.5
$ws.Range("M99").Value = 764.5
$ws.Range("H99").Value = 295.6
$ws.Range("I99").Value = 244.5
$ws.Range("I100").Value = 1319.6666
$ws.Range("H100").Value = 1241.1052
$ws.Range("K100").Value = 1319.6666
$ws.Range("M100").Value = -778.6666
$ws.Range("N112").Value = -15505.571
$ws.Range("J112").Value = 4429.857
$ws.Range("H112").Value = 3724.4546
$ws.Range("L112").Value = 13289.571
$ws.Range("L122").Value = 993.75
$ws.Range("N122").Value = -5893.75
$ws.Range("J122").Value = 331.25
$ws.Range("I122").Value = 813.3182
$ws.Range("M122").Value = 10.04539999999997
$ws.Range("H122").Value = 739.1539
$ws.Range("K122").Value = 2439.9546
$ws.Range("I125").Value = 2571468.8
$ws.Range("M125").Value = -23140759.2
$ws.Range("H125").Value = 2146224
$ws.Range("K125").Value = 23143219.2
$ws.Range("I132").Value = 1539.1364
$ws.Range("M132").Value = -2087.4092
$ws.Range("H132").Value = 1602.3334
$ws.Range("K132").Value = 4617.4092
$ws.Range("K137").Value = 6257.25
$ws.Range("M137").Value = -3707.25
$ws.Range("H137").Value = 3640.4167
$ws.Range("I137").Value = 2085.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 5208.469
$ws.Range("M32").Value = -4921.469
$ws.Range("H32").Value = 5592.9805
$ws.Range("K32").Value = 5208.469
$ws.Range("K45").Value = 127628.81
$ws.Range("M45").Value = -127251.81
$ws.Range("I45").Value = 127628.81
$ws.Range("H45").Value = 90806.78
$ws.Range("I61").Value = 10886.5
$ws.Range("M61").Value = -10674.5
$ws.Range("H61").Value = 10621.333
$ws.Range("K61").Value = 10886.5
$ws.Range("K97").Value = 2970
$ws.Range("L97").Value = 3033.3333
$ws.Range("J97").Value = 3033.3333
$ws.Range("N97").Value = -4025.3333
$ws.Range("M97").Value = -2474
$ws.Range("H97").Value = 3001.6667
$ws.Range("I97").Value = 2970
$ws.Range("M102").Value = -959.1819999999998
$ws.Range("I102").Value = 2581.182
$ws.Range("H102").Value = 3715.4119
$ws.Range("K102").Value = 2581.182
$ws.Range("L122").Value = 12823.5
$ws.Range("N122").Value = -17723.5
$ws.Range("J122").Value = 4274.5
$ws.Range("I122").Value = 2343
$ws.Range("M122").Value = -4579
$ws.Range("H122").Value = 2825.875
$ws.Range("K122").Value = 7029
$ws.Range("M136").Value = -30109.5
$ws.Range("H136").Value = 10621.333
$ws.Range("I136").Value = 10886.5
$ws.Range("K136").Value = 32659.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K94").Value = 1020.6
$ws.Range("M94").Value = -569.6
$ws.Range("H94").Value = 1020.6
$ws.Range("I94").Value = 1020.6
$ws.Range("K105").Value = 659.2273
$ws.Range("L105").Value = 983.6
$ws.Range("N105").Value = -4477.6
$ws.Range("J105").Value = 983.6
$ws.Range("I105").Value = 659.2273
$ws.Range("M105").Value = 1087.7727
$ws.Range("H105").Value = 719.2963
$ws.Range("I134").Value = 7014.8945
$ws.Range("M134").Value = -18509.6835
$ws.Range("H134").Value = 7194.381
$ws.Range("K134").Value = 21044.6835

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L9").Value = 70000
$ws.Range("J9").Value = 70000
$ws.Range("N9").Value = -70336
$ws.Range("H9").Value = 70000
$ws.Range("I132").Value = 4836.3
$ws.Range("M132").Value = -11978.9
$ws.Range("H132").Value = 7168.7144
$ws.Range("K132").Value = 14508.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K4").Value = 49742475
$ws.Range("M4").Value = -49742363
$ws.Range("H4").Value = 17360118
$ws.Range("I4").Value = 16580825
$ws.Range("L37").Value = 417735
$ws.Range("N37").Value = -417959
$ws.Range("J37").Value = 139245
$ws.Range("H37").Value = 139245
$ws.Range("J50").Value = 378.66666
$ws.Range("M50").Value = -3867.799999999999
$ws.Range("H50").Value = 865.4545000000001
$ws.Range("I50").Value = 1449.6
$ws.Range("K50").Value = 4348.799999999999
$ws.Range("L50").Value = 1135.99998
$ws.Range("N50").Value = -2097.99998
$ws.Range("K53").Value = 4348.799999999999
$ws.Range("L53").Value = 1135.99998
$ws.Range("J53").Value = 378.66666
$ws.Range("N53").Value = -2097.99998
$ws.Range("M53").Value = -3867.799999999999
$ws.Range("H53").Value = 865.4545000000001
$ws.Range("I53").Value = 1449.6
$ws.Range("H55").Value = 915.55554
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("J55").Value = 1000

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J80").Value = 4999
$ws.Range("N80").Value = -6995
$ws.Range("H80").Value = 3999.5
$ws.Range("L80").Value = 4999
$ws.Range("L83").Value = 24995
$ws.Range("N83").Value = -34979
$ws.Range("J83").Value = 4999
$ws.Range("H83").Value = 3999.5
$ws.Range("I122").Value = 2000
$ws.Range("M122").Value = -3550
$ws.Range("H122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("H126").Value = 3468.9285
$ws.Range("I126").Value = 2278.6365
$ws.Range("K126").Value = 6835.9095
$ws.Range("M126").Value = -4365.9095
$ws.Range("I132").Value = 3593.6553
$ws.Range("M132").Value = -8250.965899999999
$ws.Range("H132").Value = 3744.4375
$ws.Range("K132").Value = 10780.9659

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""
$ws.Range("J38").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("I63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("I66").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L68").Value = 3666
$ws.Range("N68").Value = -5164
$ws.Range("J68").Value = 3666
$ws.Range("I68").Value = 3565.5
$ws.Range("M68").Value = -2816.5
$ws.Range("H68").Value = 3599
$ws.Range("K68").Value = 3565.5
$ws.Range("L71").Value = 18330
$ws.Range("J71").Value = 3666
$ws.Range("N71").Value = -25818
$ws.Range("M71").Value = -14083.5
$ws.Range("H71").Value = 3599
$ws.Range("I71").Value = 3565.5
$ws.Range("K71").Value = 17827.5
$ws.Range("I132").Value = 10890.956
$ws.Range("M132").Value = -30142.868
$ws.Range("H132").Value = 10378.68
$ws.Range("K132").Value = 32672.868
$ws.Range("M136").Value = -10826.175
$ws.Range("H136").Value = 4408.488
$ws.Range("I136").Value = 4458.725
$ws.Range("K136").Value = 13376.175

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M96").Value = -415.25
$ws.Range("H96").Value = 1820.6
$ws.Range("I96").Value = 1788.25
$ws.Range("K96").Value = 1788.25
$ws.Range("L122").Value = 22499.4
$ws.Range("N122").Value = -27399.4
$ws.Range("J122").Value = 7499.8
$ws.Range("I122").Value = 2976.1538
$ws.Range("M122").Value = -6478.4614
$ws.Range("H122").Value = 4232.722
$ws.Range("K122").Value = 8928.4614
$ws.Range("I132").Value = 3879.147
$ws.Range("M132").Value = -9107.440999999999
$ws.Range("H132").Value = 4556.476
$ws.Range("K132").Value = 11637.441
$ws.Range("N138").Value = -175281
$ws.Range("H138").Value = 165001
$ws.Range("L138").Value = 165001
$ws.Range("J138").Value = 165001
